# Update cryptos list: refresh Price (column D) and Volume(1h) (column E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force the literal string into the cell even if it looks like a number,
    # then strip the resulting quote-prefix style so no style index is left on the cell.
    $range.Value = "'" + $text
    $range.Style = 'Normal'
}

$ws.Range('D2').Value = '22.395.50'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '1.570.83'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  +0.13%  '
Set-TextValue $ws.Range('D5') '1.002'
$ws.Range('E5').Value = '  +0.17%  '
Set-TextValue $ws.Range('D6') '291.75'
$ws.Range('E6').Value = '  +0.69%  '
Set-TextValue $ws.Range('D7') '0.3758'
$ws.Range('E7').Value = '  +2.21%  '
Set-TextValue $ws.Range('D8') '49.78'
$ws.Range('E8').Value = '  +0.90%  '
Set-TextValue $ws.Range('D9') '0.3414'
$ws.Range('E9').Value = '  +0.47%  '
Set-TextValue $ws.Range('D10') '0.07619'
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('E11').Value = '  -2.50%  '
$ws.Range('E12').Value = '  +0.14%  '
Set-TextValue $ws.Range('D13') '21.18'
$ws.Range('E13').Value = '  -1.25%  '
Set-TextValue $ws.Range('D14') '5.996'
$ws.Range('E14').Value = '  -1.22%  '
Set-TextValue $ws.Range('D15') '6.940'
$ws.Range('E15').Value = '  -0.02%  '
$ws.Range('D16').Value = '1.569.53'
$ws.Range('E16').Value = '  +0.24%  '
Set-TextValue $ws.Range('D17') '0.00001134'
$ws.Range('E17').Value = '  -0.57%  '
Set-TextValue $ws.Range('D18') '90.25'
$ws.Range('E18').Value = '  +0.94%  '
Set-TextValue $ws.Range('D19') '0.06747'
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('E20').Value = '  +0.12%  '
Set-TextValue $ws.Range('D21') '16.76'
$ws.Range('E21').Value = '  +0.95%  '
Set-TextValue $ws.Range('D22') '6.189'
$ws.Range('E22').Value = '  -1.10%  '
Set-TextValue $ws.Range('D23') '11.99'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').Value = '22.388.40'
$ws.Range('E24').Value = '  -0.14%  '
Set-TextValue $ws.Range('D25') '2.393'
$ws.Range('E25').Value = '  +0.50%  '
Set-TextValue $ws.Range('D26') '2.676'
$ws.Range('E26').Value = '  -9.69%  '
Set-TextValue $ws.Range('D27') '20.12'
$ws.Range('E27').Value = '  +0.53%  '
Set-TextValue $ws.Range('D28') '147.25'
$ws.Range('E28').Value = '  +0.68%  '
Set-TextValue $ws.Range('D29') '5.048'
$ws.Range('E29').Value = '  +1.36%  '
Set-TextValue $ws.Range('D30') '126.61'
$ws.Range('E30').Value = '  +0.49%  '
$ws.Range('D31').Value = '1.747.04'
$ws.Range('E31').Value = '  +0.31%  '
Set-TextValue $ws.Range('D32') '2.013'
$ws.Range('E32').Value = '  +0.37%  '
Set-TextValue $ws.Range('D33') '6.096'
$ws.Range('E33').Value = '  -3.31%  '
$ws.Range('E34').Value = '  -5.27%  '
Set-TextValue $ws.Range('D35') '10.12'
$ws.Range('E35').Value = '  -1.95%  '
Set-TextValue $ws.Range('D36') '0.08502'
$ws.Range('E36').Value = '  +0.45%  '
Set-TextValue $ws.Range('D37') '0.02534'
$ws.Range('E37').Value = '  -0.41%  '
Set-TextValue $ws.Range('D38') '1.378'
$ws.Range('E38').Value = '  +10.38%  '
$ws.Range('E39').Value = '  -1.20%  '
Set-TextValue $ws.Range('D40') '0.06495'
$ws.Range('E40').Value = '  -1.16%  '
Set-TextValue $ws.Range('D41') '5.415'
$ws.Range('E41').Value = '  -2.86%  '
$ws.Range('E42').Value = '  -3.43%  '
Set-TextValue $ws.Range('D43') '0.6345'
$ws.Range('E43').Value = '  -0.78%  '
Set-TextValue $ws.Range('D44') '1.001'
$ws.Range('E44').Value = '  +0.07%  '
Set-TextValue $ws.Range('D45') '14.02'
$ws.Range('E45').Value = '  -2.14%  '
Set-TextValue $ws.Range('D46') '3.794'
$ws.Range('E46').Value = '  +1.30%  '
Set-TextValue $ws.Range('D47') '0.5948'
$ws.Range('E47').Value = '  -1.11%  '
Set-TextValue $ws.Range('D48') '2.087'
$ws.Range('E48').Value = '  -1.82%  '
$ws.Range('E49').Value = '  +0.99%  '
Set-TextValue $ws.Range('D50') '123.98'
$ws.Range('E50').Value = '  +0.28%  '
Set-TextValue $ws.Range('D51') '0.07322'
$ws.Range('E51').Value = '  +0.26%  '
